$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Majors"

$ws.Range("A1").Value = "Kode Jurusan"
$ws.Range("B1").Value = "Jurusan"
$ws.Range("C1").Value = "Jenjang"

$ws.Columns.Item(1).ColumnWidth = 15

$ws.Range("C2").Select()
